# Minor label edits: the "risk factors" summary rows on the "data" and
# "languages" sheets were not age-scoped; clarify that they refer to the
# population aged 18-69 years (English + French wording).

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")
$wsLang = $wb.Worksheets.Item("languages")

# --- "data" sheet: column C holds the English label text ---------------
$wsData.Range("C103").Value = "of the population aged 18-69 years have 3 to 5 risk factors"
$wsData.Range("C104").Value = "of the population aged 18-69 years have 1 to 2 risk factors"
$wsData.Range("C105").Value = "of the population aged 18-69 years have 0 risk factors"

# --- "languages" sheet: column A = English, column B = French ----------
$wsLang.Range("A103").Value = "of the population aged 18-69 years have 3 to 5 risk factors"
$wsLang.Range("B103").Value = "de la population âgées de 18 à 69 ans présentent 3 à 5 facteurs de risque"

$wsLang.Range("A104").Value = "of the population aged 18-69 years have 1 to 2 risk factors"
$wsLang.Range("B104").Value = "de la population âgées de 18 à 69 ans présentent 1 à 2 facteurs de risque"

$wsLang.Range("A105").Value = "of the population aged 18-69 years have 0 risk factors"
$wsLang.Range("B105").Value = "de la population âgées de 18 à 69 ans ne présente aucun facteur de risque"

# --- Restore the view state left behind by the edit ---------------------
# The "data" sheet became the active tab, scrolled down to the edited rows.
$wsData.Activate()
$excel.ActiveWindow.ScrollRow = 96
$wsData.Range("C108").Select()

# The "languages" sheet keeps the edited rows selected, scrolled into view.
$wsLang.Activate()
$excel.ActiveWindow.ScrollRow = 91
$wsLang.Range("A103:A105").Select()

$wsData.Activate()
